$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.966.30"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.555.35"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D5").Value = "'207.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").Value = "'22.07"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.13%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "1.776.02"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.556.83"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "26.956.00"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'61.76"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").Value = "'217.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "'7.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").Value = "'153.30"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'15.01"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'0.0470"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").Value = "'3.11"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.98%  "
$ws.Range("D34").Value = "1.418.70"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("E36").Value = "  +11.94%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").Value = "'0.529"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").Value = "'0.807"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "'64.53"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").Value = "'1.74"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "1.689.52"
$ws.Range("D48").Value = "'87.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").Value = "0.0₇0989"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("E51").Value = "  +0.75%  "
